# Append a new row (row 6) to the Results sheet, mirroring the existing
# rows but representing a newly-finished "sat" run for group 6-5-5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "6-5-5"
$ws.Range("C6").Value = "new_sequential"

# D6 must stay text ("78.324"), not be coerced into a numeric value.
# Force text storage the same way real Excel COM automation does:
# mark the cell as Text before assigning, then drop back to the
# "Normal" style so no stray number format sticks to the cell itself.
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "78.324"
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = "sat"
$ws.Range("F6").Value = 4140
$ws.Range("G6").Value = 171292
